$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting the existing rows 30-32 down to 31-33.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new weekly data point.
$ws.Range("A30").Value = 5
$ws.Range("B30").Value = "Macroferia Regional de Talca"
$ws.Range("C30").Value = "Maule"
$ws.Range("D30").Value = 44476
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = 100112026
$ws.Range("G30").Value = "Haba"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 500
$ws.Range("K30").Value = 8000
$ws.Range("L30").Value = 8000
$ws.Range("M30").Value = 8000
$ws.Range("N30").Value = "$/saco 25 kilos"
$ws.Range("O30").Value = "Región de O'Higgins"
$ws.Range("P30").Value = 320
$ws.Range("Q30").Value = 25
$ws.Range("R30").Value = "Hortaliza"
